$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'251.70"
$ws.Range("D3").Formula = "'22.77"
$ws.Range("D4").Formula = "'5.431"
$ws.Range("D5").Formula = "'0.05692"
$ws.Range("D6").Formula = "'3.411"
$ws.Range("D7").Formula = "'6.386"
$ws.Range("D8").Formula = "'0.8129"
$ws.Range("D9").Formula = "'0.9413"
$ws.Range("D10").Formula = "'0.1436"
$ws.Range("D11").Formula = "'0.07491"
$ws.Range("D12").Formula = "'0.03161"
$ws.Range("D13").Formula = "'0.03079"
$ws.Range("D14").Formula = "'0.09363"
$ws.Range("D15").Formula = "'3.731"
$ws.Range("D16").Formula = "'0.001604"
$ws.Range("D17").Formula = "'0.04762"
$ws.Range("D18").Formula = "'0.0005787"
$ws.Range("D19").Formula = "'0.006372"
$ws.Range("D20").Formula = "'0.005045"
$ws.Range("D21").Formula = "'0.001029"
$ws.Range("D22").Formula = "'0.0001500"
$ws.Range("D23").Formula = "'3.710"
$ws.Range("D24").Formula = "'2.170"
$ws.Range("D26").Formula = "'0.1307"
$ws.Range("D28").Formula = "'0.0002998"
$ws.Range("D41").Formula = "'0.006798"
$ws.Range("D42").Formula = "'0.1072"
$ws.Range("D43").Formula = "'0.002709"
$ws.Range("D44").Formula = "'0.008125"
$ws.Range("D45").Formula = "'0.00005757"
$ws.Range("D47").Formula = "'0.4997"
$ws.Range("D49").Formula = "'0.00002099"
$ws.Range("D50").Formula = "'0.01009"
